# "Version 2." -> "Version 1." while collapsing the run structure to
# match the target: "Versi"+"on" merge into a single "Version" run, the
# " 2" run becomes " 1.", and the trailing "." run (after the bookmark)
# is removed entirely.

$d = $word.ActiveDocument

# Step 1: turn " 2" into " 1." (this run sits between spellEnd and the
# bookmark, so its boundaries are untouched).
$rNum = $d.Range(7, 9)
$rNum.Text = " 1."

# Step 2: merge the "Versi" / "on" runs into a single "Version" run.
# A direct same-text assignment is a no-op for the engine's run-diffing,
# so nudge it through an intermediate value first to force the merge.
$rWord = $d.Range(0, 7)
$rWord.Text = "VersionX"
$rWord2 = $d.Range(0, 8)
$rWord2.Text = "Version"

# Step 3: delete the now-orphaned trailing "." run that followed the
# bookmark (its text got folded into " 1." above).
$rTrail = $d.Range(10, 11)
$rTrail.Delete()
